# documentation/simulation_parameters.xlsx
#
# Swap out the old "constantClimate ... resetSWC" parameter table for the
# new SurEau-Ecos simulation-setup table (mainDir ... overWrite). The new
# table is 7 parameters (+ header) instead of 10, so the sheet shrinks
# from B2:C12 down to B2:C9; the two rows that used to carry the yellow
# "not-yet-implemented" highlight fill lose it, and the bottom border that
# used to close the table at row 12 now closes it at row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Overwrite rows 3-9 in place, so every cell keeps the row/column
#    base style (Courier-new param names, wrapped Arial descriptions)
#    it already had - only the two content-driven tweaks below need to
#    touch formatting at all.
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "mainDir"
$ws.Range("C3").Value = "the main directory of the model "

$ws.Range("B4").Value = "startYearSimulation"
$ws.Range("C4").Value = "a numeric indicating the starting year for the simulation (must match the dates of the input climate file)"

$ws.Range("B5").Value = "endYearSimulation"
$ws.Range("C5").Value = "a numeric indicating the last year for the simulation (must match the dates of the input climate file)"

$ws.Range("B6").Value = "resolutionOutput"
$ws.Range("C6").Value = "the resolution for the output simulation file. Must be 'subdaily', 'daily' or 'yearly'"

$ws.Range("B7").Value = "outputType"
$ws.Range("C7").Value = "the type of output. "

$ws.Range("B8").Value = "outputPath"
$ws.Range("C8").Value = "the path of the ouput simulation fire  "

$ws.Range("B9").Value = "overWrite"
$ws.Range("C9").Value = "a logical value indicating wether the chosen output path can be overwrittent if it already exists (default = F)"

# ---------------------------------------------------------------------
# 2. Drop the old yellow "not implemented" highlight on the two rows
#    that used to carry it (previously compOptionsForEvapo /
#    numericalScheme) - resolutionOutput's description keeps the
#    quote-prefixed flavour of that second one.
# ---------------------------------------------------------------------
$ws.Range("C5").Interior.Pattern = -4142
$ws.Range("C6").Interior.Pattern = -4142

# ---------------------------------------------------------------------
# 3. Row heights: the two-line entries need 28pt, everything else goes
#    back to the sheet's default (AutoFit removes any leftover explicit
#    height from the old layout).
# ---------------------------------------------------------------------
$ws.Rows("3:3").AutoFit()
$ws.Rows("4:4").RowHeight = 28
$ws.Rows("5:5").RowHeight = 28
$ws.Rows("6:6").RowHeight = 28
$ws.Rows("7:7").AutoFit()
$ws.Rows("8:8").AutoFit()
$ws.Rows("9:9").RowHeight = 28

# ---------------------------------------------------------------------
# 4. The table-closing bottom border moves from row 12 to row 9 - copy
#    it across (re-uses the existing border style instead of minting a
#    new one) before the trailing rows are removed.
# ---------------------------------------------------------------------
$ws.Range("B12:C12").Copy()
$ws.Range("B9:C9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5. The old rows 10-12 (avoidWaterSoilTransfer, defoliation, resetSWC)
#    are gone in the new table - delete them so the sheet shrinks back
#    down to B2:C9.
# ---------------------------------------------------------------------
$ws.Rows("10:12").Delete()

# ---------------------------------------------------------------------
# 6. Restore the view state recorded for the edited sheet.
# ---------------------------------------------------------------------
$win = $ws.Application.ActiveWindow
$win.Zoom = 150
$win.ScrollRow = 2
$ws.Range("C12").Select()
